# Applies the "rajoue de add eparne est des gain et depence" commit:
# clears out the sample/demo transaction rows on the four "flux" sheets
# (deponce_continue, deponce_courante, revenu_continue, revenu_courante)
# while leaving the date-styled placeholder cells in column A (and B for
# the "courante" sheets) in place, and repoints the saved selection /
# active-tab state at the "eparne" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# deponce_continue: clear the three demo expense rows (netflix/Crunchyrol/ADN)
# ---------------------------------------------------------------------
$wsDeponceContinue = $wb.Worksheets.Item("deponce_continue")
$wsDeponceContinue.Range("B2:C2").ClearContents()
$wsDeponceContinue.Range("A2").ClearContents()
$wsDeponceContinue.Range("B3:C3").ClearContents()
$wsDeponceContinue.Range("A3").ClearContents()
$wsDeponceContinue.Range("B4:C4").ClearContents()
$wsDeponceContinue.Range("A4").ClearContents()

# ---------------------------------------------------------------------
# deponce_courante: clear the demo "cours" expense row
# ---------------------------------------------------------------------
$wsDeponceCourante = $wb.Worksheets.Item("deponce_courante")
$wsDeponceCourante.Range("C2:D2").ClearContents()
$wsDeponceCourante.Range("A2:B2").ClearContents()

# ---------------------------------------------------------------------
# revenu_continue: clear the demo "salaire" income row
# ---------------------------------------------------------------------
$wsRevenuContinue = $wb.Worksheets.Item("revenu_continue")
$wsRevenuContinue.Range("B2:C2").ClearContents()
$wsRevenuContinue.Range("A2").ClearContents()

# ---------------------------------------------------------------------
# revenu_courante: clear the demo "prime" income row
# ---------------------------------------------------------------------
$wsRevenuCourante = $wb.Worksheets.Item("revenu_courante")
$wsRevenuCourante.Range("C2:D2").ClearContents()
$wsRevenuCourante.Range("A2:B2").ClearContents()

# ---------------------------------------------------------------------
# Restore each sheet's remembered selection, then finish on "eparne" so
# it becomes the workbook's active tab/sheet again.
# ---------------------------------------------------------------------
$wsDeponceContinue.Activate() | Out-Null
$wsDeponceContinue.Range("A4").Select() | Out-Null

$wsDeponceCourante.Activate() | Out-Null
$wsDeponceCourante.Range("D2").Select() | Out-Null

$wsRevenuContinue.Activate() | Out-Null
$wsRevenuContinue.Range("C2").Select() | Out-Null

$wsRevenuCourante.Activate() | Out-Null
$wsRevenuCourante.Range("C2").Select() | Out-Null

$wsEparne = $wb.Worksheets.Item("eparne")
$wsEparne.Activate() | Out-Null
$wsEparne.Range("D6").Select() | Out-Null
